$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small numeric corrections (rounding adjustments) ---
$ws.Range("E57").Value = 82217.992
$ws.Range("I57").Value = 108581
$ws.Range("M57").Value = 110746.992

$ws.Range("I59").Value = 90250.008
$ws.Range("Q59").Value = 81705
$ws.Range("U59").Value = 83035
$ws.Range("Y59").Value = 87228.016
$ws.Range("AG59").Value = 114814

$ws.Range("Q60").Value = -49107
$ws.Range("Y60").Value = -54099.992
$ws.Range("AG60").Value = -66225.992

# --- Clear stray placeholder zeros (now blank, concatenated balance sheets) ---
$ws.Range("Q57:AJ57").ClearContents()
$ws.Range("Q58:AJ58").ClearContents()
$ws.Range("B64:Q64").ClearContents()
$ws.Range("Q71:AJ71").ClearContents()
$ws.Range("Q72:AJ72").ClearContents()
$ws.Range("Q73:AJ73").ClearContents()
$ws.Range("Q77:AJ77").ClearContents()
$ws.Range("Q78:AJ78").ClearContents()
